$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: merge the two runs "THU Oct 19" + " 11:06:47 IST 2018" into a
# single run "THU Oct 19 11:06:47 IST 2018". Both runs already share the
# same formatting, so a Find/Replace over the already-matching text
# coalesces the text into one run (mirrors what Word does when you
# select-and-retype text that already reads correctly).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("THU Oct 19 11:06:47 IST 2018", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "THU Oct 19 11:06:47 IST 2018", 2) | Out-Null

# ---------------------------------------------------------------------
# Part 2: locate the "Amount Received mode ... - CASH AND CLEARD" line
# that closes out the very last purchase entry, and splice in a brand
# new purchase block (MAMATHA CHICK IN, 25/10/2018) right after it -
# before the pre-existing trailing blank paragraphs.
# ---------------------------------------------------------------------
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*CASH AND CLEARD*") {
        $anchorIndex = $i
    }
}

$anchor = $d.Paragraphs($anchorIndex).Range

# Pre-create 13 blank paragraphs right after the anchor. Re-using the
# same (fixed) anchor range for every InsertParagraphAfter() call makes
# each new paragraph land immediately after it, in call order, so the
# resulting indices anchorIndex+1 .. anchorIndex+13 come out in the
# correct forward order.
for ($k = 0; $k -lt 13; $k++) {
    $anchor.InsertParagraphAfter()
}

$base = $anchorIndex

# base+1 : blank separator line (already blank - nothing to do)

# base+2 : date line, split into two runs "TUE Oct 23" + " 14:05:08 IST 2018"
$p = $d.Paragraphs($base + 2)
$p.Range.InsertAfter("TUE Oct 23 14:05:08 IST 2018")
$lineStart = $p.Range.Start
$splitAt = $lineStart + ("TUE Oct 23").Length
$lineEnd = $p.Range.End
$sub = $d.Range($splitAt, $lineEnd)
$sub.Font.Bold = 1
$sub.Font.Bold = 0

# base+3 : Person Name
$d.Paragraphs($base + 3).Range.InsertAfter("Person Name`t`t`t`t- HANUMANTHARAYA")

# base+4 : Bill number
$d.Paragraphs($base + 4).Range.InsertAfter("Bill number`t`t`t`t- 8563")

# base+5 : divider
$d.Paragraphs($base + 5).Range.InsertAfter("---------------------------------------------------------------")

# base+6 : Item Name
$d.Paragraphs($base + 6).Range.InsertAfter("Item Name`t`t`t`t- CHOWCHOW")

# base+7 : Number of Pockets
$d.Paragraphs($base + 7).Range.InsertAfter("Number of Pockets`t`t`t- 1")

# base+8 : Number of KGs
$d.Paragraphs($base + 8).Range.InsertAfter("Number of KGs`t`t`t- 67")

# base+9 : Rate
$d.Paragraphs($base + 9).Range.InsertAfter("Rate`t`t`t`t`t- 7")

# base+10 : Total Price
$d.Paragraphs($base + 10).Range.InsertAfter("Total Price`t`t`t`t- 469.0")

# base+11 : Amount balance (bold)
$p11 = $d.Paragraphs($base + 11)
$p11.Range.Font.Bold = 1
$p11.Range.InsertAfter("Amount balance`t`t`t- 469.0")

# base+12, base+13 : trailing blank lines (already blank - nothing to do)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
